# Weekly update: insert a new price record (new reporting date) as row 5,
# pushing the existing historical rows (old rows 5-47) down by one to rows 6-48.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 5 (shifts rows 5..47 down to 6..48,
# carrying their formatting/styles along, and extends the sheet dimension).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with this week's record.
$ws.Cells.Item(5, 1).Value = 8
$ws.Cells.Item(5, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(5, 3).Value = "Coquimbo"
$ws.Cells.Item(5, 4).Value = 44490
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 100112052
$ws.Cells.Item(5, 7).Value = "Albahaca"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 660
$ws.Cells.Item(5, 11).Value = 3500
$ws.Cells.Item(5, 12).Value = 4000
$ws.Cells.Item(5, 13).Value = 3750
$ws.Cells.Item(5, 14).Value = "$/paquete"
$ws.Cells.Item(5, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 16).Value = 3750
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
